$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.622599999999997

$ws.Range("B4").Value = 4.759800000000005
$ws.Range("C4").Value = -14.07400000000001
$ws.Range("D4").Value = -8.007400000000001

$ws.Range("C5").Value = -14.77670000000001

$ws.Range("B6").Value = 9.316799999999995

$ws.Range("B7").Value = 6.0039

$ws.Range("C8").Value = -12.0485

$ws.Range("D9").Value = -7.865699999999999

$ws.Range("D11").Value = -8.375400000000001

$ws.Range("D14").Value = -6.550699999999997

$ws.Range("B16").Value = 9.105400000000005
$ws.Range("C16").Value = -12.3814

$ws.Range("D18").Value = -8.528199999999995

$ws.Range("B20").Value = 5.499099999999998

$ws.Range("C22").Value = -10.97349999999999

$ws.Range("D25").Value = -8.462899999999992
